$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.023.89"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "2.269.37"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("B5").Value = "XRP"
$ws.Range("C5").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.655"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +4.51%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "233.08"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +0.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.92"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.453"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +6.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0982"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.11"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.89"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +3.53%  "
$ws.Range("E13").Value = "  +1.87%  "
$ws.Range("D14").Value = "2.602.36"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("E16").Value = "  +4.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.842"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +3.20%  "
$ws.Range("D18").Value = "2.265.32"
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("D19").Value = "43.947.93"
$ws.Range("E19").Value = "  +0.68%  "
$ws.Range("E20").Value = "  +3.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "74.19"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.02"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.45"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -2.61%  "
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.32"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +18.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.98"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.25"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +8.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "173.79"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.127"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +4.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.05"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +7.20%  "
$ws.Range("E35").Value = "  -0.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.01"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -1.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.70"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -2.93%  "
$ws.Range("E38").Value = "  -4.73%  "
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("E40").Value = "  +3.47%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.83"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +4.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.000222"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.40"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +3.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "98.65"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0953"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -1.36%  "
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.39"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -3.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.37"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +2.60%  "
$ws.Range("D50").Value = "1.453.17"
$ws.Range("E50").Value = "  -1.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.98"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -7.06%  "
